# This script reproduces the weekly "Fruta / hortaliza" update for the
# Platano sheet: the block of 3 rows (Pinton / Primera Maduro / Primera
# Pinton) for each reporting date in the A1:T586 range is shifted down by
# one block (i.e. a new, more recent, reporting date is inserted at the top
# of the range), and the oldest block that falls off the bottom of the
# original range is appended as 3 brand-new rows, growing the sheet from
# A1:T586 to A1:T589.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that never vary across this block of rows.
$static = @{A=8; B='Terminal La Palmera de La Serena'; C='Coquimbo'; E=4; F='Fruta';
            G=100108; H='Tropicales y subtropicales'; I=100108006; J='Plátano';
            K='Sin especificar'; Q='$/caja 20 kilos'; R='Ecuador'; T=20}

# Target (date, quality, volume, min/max/avg price, price-per-kg) for every
# row from 476 to 589, post-shift.
$rowsData = @(
    @{R=476;D=44642;L='Pintón';M=120;N=18000;O=18000;P=18000;S=900},
    @{R=477;D=44642;L='Primera Maduro';M=200;N=20000;O=20000;P=20000;S=1000},
    @{R=478;D=44642;L='Primera Pintón';M=200;N=21000;O=21000;P=21000;S=1050},
    @{R=479;D=44473;L='Pintón';M=80;N=21000;O=21000;P=21000;S=1050},
    @{R=480;D=44473;L='Primera Maduro';M=120;N=22000;O=22000;P=22000;S=1100},
    @{R=481;D=44473;L='Primera Pintón';M=120;N=23000;O=23000;P=23000;S=1150},
    @{R=482;D=44357;L='Pintón';M=80;N=13000;O=13000;P=13000;S=650},
    @{R=483;D=44357;L='Primera Maduro';M=120;N=14500;O=14500;P=14500;S=725},
    @{R=484;D=44357;L='Primera Pintón';M=120;N=15000;O=15000;P=15000;S=750},
    @{R=485;D=44455;L='Pintón';M=80;N=18000;O=18000;P=18000;S=900},
    @{R=486;D=44455;L='Primera Maduro';M=120;N=19000;O=19000;P=19000;S=950},
    @{R=487;D=44455;L='Primera Pintón';M=120;N=20000;O=20000;P=20000;S=1000},
    @{R=488;D=44581;L='Pintón';M=80;N=15000;O=15000;P=15000;S=750},
    @{R=489;D=44581;L='Primera Maduro';M=120;N=17000;O=17000;P=17000;S=850},
    @{R=490;D=44581;L='Primera Pintón';M=120;N=18000;O=18000;P=18000;S=900},
    @{R=491;D=44553;L='Pintón';M=120;N=12000;O=12000;P=12000;S=600},
    @{R=492;D=44553;L='Primera Maduro';M=200;N=14000;O=14000;P=14000;S=700},
    @{R=493;D=44553;L='Primera Pintón';M=200;N=15000;O=15000;P=15000;S=750},
    @{R=494;D=44462;L='Pintón';M=80;N=16000;O=16000;P=16000;S=800},
    @{R=495;D=44462;L='Primera Maduro';M=120;N=17000;O=17000;P=17000;S=850},
    @{R=496;D=44462;L='Primera Pintón';M=120;N=18000;O=18000;P=18000;S=900},
    @{R=497;D=44641;L='Pintón';M=80;N=18000;O=18000;P=18000;S=900},
    @{R=498;D=44641;L='Primera Maduro';M=120;N=20000;O=20000;P=20000;S=1000},
    @{R=499;D=44641;L='Primera Pintón';M=120;N=21000;O=21000;P=21000;S=1050},
    @{R=500;D=44421;L='Pintón';M=80;N=14000;O=14000;P=14000;S=700},
    @{R=501;D=44421;L='Primera Maduro';M=120;N=15500;O=15500;P=15500;S=775},
    @{R=502;D=44421;L='Primera Pintón';M=120;N=16000;O=16000;P=16000;S=800},
    @{R=503;D=44329;L='Pintón';M=80;N=15000;O=15000;P=15000;S=750},
    @{R=504;D=44329;L='Primera Maduro';M=120;N=16500;O=16500;P=16500;S=825},
    @{R=505;D=44329;L='Primera Pintón';M=120;N=17000;O=17000;P=17000;S=850},
    @{R=506;D=44637;L='Pintón';M=80;N=18000;O=18000;P=18000;S=900},
    @{R=507;D=44637;L='Primera Maduro';M=120;N=20000;O=20000;P=20000;S=1000},
    @{R=508;D=44637;L='Primera Pintón';M=120;N=21000;O=21000;P=21000;S=1050},
    @{R=509;D=44445;L='Pintón';M=80;N=19000;O=19000;P=19000;S=950},
    @{R=510;D=44445;L='Primera Maduro';M=120;N=20000;O=20000;P=20000;S=1000},
    @{R=511;D=44445;L='Primera Pintón';M=120;N=21000;O=21000;P=21000;S=1050},
    @{R=512;D=44530;L='Pintón';M=120;N=21000;O=21000;P=21000;S=1050},
    @{R=513;D=44530;L='Primera Maduro';M=120;N=22000;O=22000;P=22000;S=1100},
    @{R=514;D=44530;L='Primera Pintón';M=120;N=23000;O=23000;P=23000;S=1150},
    @{R=515;D=44483;L='Pintón';M=80;N=22000;O=22000;P=22000;S=1100},
    @{R=516;D=44483;L='Primera Maduro';M=120;N=23000;O=23000;P=23000;S=1150},
    @{R=517;D=44483;L='Primera Pintón';M=120;N=24000;O=24000;P=24000;S=1200},
    @{R=518;D=44294;L='Pintón';M=80;N=15000;O=15000;P=15000;S=750},
    @{R=519;D=44294;L='Primera Maduro';M=120;N=16500;O=16500;P=16500;S=825},
    @{R=520;D=44294;L='Primera Pintón';M=120;N=17000;O=17000;P=17000;S=850},
    @{R=521;D=44617;L='Pintón';M=80;N=16000;O=16000;P=16000;S=800},
    @{R=522;D=44617;L='Primera Maduro';M=120;N=18000;O=18000;P=18000;S=900},
    @{R=523;D=44617;L='Primera Pintón';M=120;N=19000;O=19000;P=19000;S=950},
    @{R=524;D=44557;L='Pintón';M=120;N=12000;O=12000;P=12000;S=600},
    @{R=525;D=44557;L='Primera Maduro';M=160;N=14000;O=14000;P=14000;S=700},
    @{R=526;D=44557;L='Primera Pintón';M=160;N=15000;O=15000;P=15000;S=750},
    @{R=527;D=44396;L='Pintón';M=80;N=14000;O=14000;P=14000;S=700},
    @{R=528;D=44396;L='Primera Maduro';M=120;N=15500;O=15500;P=15500;S=775},
    @{R=529;D=44396;L='Primera Pintón';M=120;N=16000;O=16000;P=16000;S=800},
    @{R=530;D=44232;L='Pintón';M=80;N=16000;O=16000;P=16000;S=800},
    @{R=531;D=44232;L='Primera Maduro';M=120;N=17500;O=17500;P=17500;S=875},
    @{R=532;D=44232;L='Primera Pintón';M=120;N=18000;O=18000;P=18000;S=900},
    @{R=533;D=44330;L='Pintón';M=80;N=15000;O=15000;P=15000;S=750},
    @{R=534;D=44330;L='Primera Maduro';M=120;N=16500;O=16500;P=16500;S=825},
    @{R=535;D=44330;L='Primera Pintón';M=120;N=17000;O=17000;P=17000;S=850},
    @{R=536;D=44504;L='Pintón';M=80;N=18000;O=18000;P=18000;S=900},
    @{R=537;D=44504;L='Primera Maduro';M=120;N=19000;O=19000;P=19000;S=950},
    @{R=538;D=44504;L='Primera Pintón';M=120;N=20000;O=20000;P=20000;S=1000},
    @{R=539;D=44572;L='Pintón';M=120;N=15000;O=15000;P=15000;S=750},
    @{R=540;D=44572;L='Primera Maduro';M=160;N=17000;O=17000;P=17000;S=850},
    @{R=541;D=44572;L='Primera Pintón';M=160;N=18000;O=18000;P=18000;S=900},
    @{R=542;D=44301;L='Pintón';M=80;N=14000;O=14000;P=14000;S=700},
    @{R=543;D=44301;L='Primera Maduro';M=120;N=15500;O=15500;P=15500;S=775},
    @{R=544;D=44301;L='Primera Pintón';M=120;N=16000;O=16000;P=16000;S=800},
    @{R=545;D=44174;L='Pintón';M=80;N=16500;O=16500;P=16500;S=825},
    @{R=546;D=44174;L='Primera Maduro';M=120;N=18000;O=18000;P=18000;S=900},
    @{R=547;D=44174;L='Primera Pintón';M=120;N=18500;O=18500;P=18500;S=925},
    @{R=548;D=44200;L='Pintón';M=80;N=17000;O=17000;P=17000;S=850},
    @{R=549;D=44200;L='Primera Maduro';M=120;N=18500;O=18500;P=18500;S=925},
    @{R=550;D=44200;L='Primera Pintón';M=120;N=19000;O=19000;P=19000;S=950},
    @{R=551;D=44385;L='Pintón';M=80;N=12000;O=12000;P=12000;S=600},
    @{R=552;D=44385;L='Primera Maduro';M=120;N=13500;O=13500;P=13500;S=675},
    @{R=553;D=44385;L='Primera Pintón';M=120;N=14000;O=14000;P=14000;S=700},
    @{R=554;D=44221;L='Pintón';M=80;N=16000;O=16000;P=16000;S=800},
    @{R=555;D=44221;L='Primera Maduro';M=120;N=17500;O=17500;P=17500;S=875},
    @{R=556;D=44221;L='Primera Pintón';M=120;N=18000;O=18000;P=18000;S=900},
    @{R=557;D=44413;L='Pintón';M=80;N=14000;O=14000;P=14000;S=700},
    @{R=558;D=44413;L='Primera Maduro';M=120;N=15500;O=15500;P=15500;S=775},
    @{R=559;D=44413;L='Primera Pintón';M=120;N=16000;O=16000;P=16000;S=800},
    @{R=560;D=44214;L='Pintón';M=80;N=14000;O=14000;P=14000;S=700},
    @{R=561;D=44214;L='Primera Maduro';M=120;N=15500;O=15500;P=15500;S=775},
    @{R=562;D=44214;L='Primera Pintón';M=120;N=16000;O=16000;P=16000;S=800},
    @{R=563;D=44610;L='Pintón';M=80;N=15000;O=15000;P=15000;S=750},
    @{R=564;D=44610;L='Primera Maduro';M=120;N=17000;O=17000;P=17000;S=850},
    @{R=565;D=44610;L='Primera Pintón';M=120;N=18000;O=18000;P=18000;S=900},
    @{R=566;D=44312;L='Pintón';M=80;N=13000;O=13000;P=13000;S=650},
    @{R=567;D=44312;L='Primera Maduro';M=120;N=14500;O=14500;P=14500;S=725},
    @{R=568;D=44312;L='Primera Pintón';M=120;N=15000;O=15000;P=15000;S=750},
    @{R=569;D=44399;L='Pintón';M=80;N=17000;O=17000;P=17000;S=850},
    @{R=570;D=44399;L='Primera Maduro';M=80;N=18500;O=18500;P=18500;S=925},
    @{R=571;D=44399;L='Primera Pintón';M=120;N=19000;O=19000;P=19000;S=950},
    @{R=572;D=44522;L='Pintón';M=80;N=22000;O=22000;P=22000;S=1100},
    @{R=573;D=44522;L='Primera Maduro';M=120;N=23000;O=23000;P=23000;S=1150},
    @{R=574;D=44522;L='Primera Pintón';M=120;N=24000;O=24000;P=24000;S=1200},
    @{R=575;D=44543;L='Pintón';M=80;N=14000;O=14000;P=14000;S=700},
    @{R=576;D=44543;L='Primera Maduro';M=120;N=16000;O=16000;P=16000;S=800},
    @{R=577;D=44543;L='Primera Pintón';M=120;N=17000;O=17000;P=17000;S=850},
    @{R=578;D=44277;L='Pintón';M=80;N=14000;O=14000;P=14000;S=700},
    @{R=579;D=44277;L='Primera Maduro';M=120;N=15500;O=15500;P=15500;S=775},
    @{R=580;D=44277;L='Primera Pintón';M=120;N=16000;O=16000;P=16000;S=800},
    @{R=581;D=44498;L='Pintón';M=80;N=21000;O=21000;P=21000;S=1050},
    @{R=582;D=44498;L='Primera Maduro';M=120;N=22500;O=22500;P=22500;S=1125},
    @{R=583;D=44498;L='Primera Pintón';M=120;N=23000;O=23000;P=23000;S=1150},
    @{R=584;D=44179;L='Pintón';M=80;N=14000;O=14000;P=14000;S=700},
    @{R=585;D=44179;L='Primera Maduro';M=120;N=15500;O=15500;P=15500;S=775},
    @{R=586;D=44179;L='Primera Pintón';M=120;N=16000;O=16000;P=16000;S=800},
    @{R=587;D=44595;L='Pintón';M=80;N=14000;O=14000;P=14000;S=700},
    @{R=588;D=44595;L='Primera Maduro';M=120;N=16000;O=16000;P=16000;S=800},
    @{R=589;D=44595;L='Primera Pintón';M=120;N=17000;O=17000;P=17000;S=850}
)

$lastExistingRow = 586

foreach ($item in $rowsData) {
    $r = $item.R
    if ($r -gt $lastExistingRow) {
        # Brand-new row: populate every column.
        $ws.Cells.Item($r, 1).Value  = $static.A
        $ws.Cells.Item($r, 2).Value  = $static.B
        $ws.Cells.Item($r, 3).Value  = $static.C
        $ws.Cells.Item($r, 5).Value  = $static.E
        $ws.Cells.Item($r, 6).Value  = $static.F
        $ws.Cells.Item($r, 7).Value  = $static.G
        $ws.Cells.Item($r, 8).Value  = $static.H
        $ws.Cells.Item($r, 9).Value  = $static.I
        $ws.Cells.Item($r, 10).Value = $static.J
        $ws.Cells.Item($r, 11).Value = $static.K
        $ws.Cells.Item($r, 17).Value = $static.Q
        $ws.Cells.Item($r, 18).Value = $static.R
        $ws.Cells.Item($r, 20).Value = $static.T
        $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }

    $ws.Cells.Item($r, 4).Value  = $item.D
    $ws.Cells.Item($r, 12).Value = $item.L
    $ws.Cells.Item($r, 13).Value = $item.M
    $ws.Cells.Item($r, 14).Value = $item.N
    $ws.Cells.Item($r, 15).Value = $item.O
    $ws.Cells.Item($r, 16).Value = $item.P
    $ws.Cells.Item($r, 19).Value = $item.S
}

